$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to text format first so numeric-looking strings
# (e.g. "1.010", "6.070") keep their exact text representation instead of
# being auto-converted into numbers by Excel.
$dCells = @("D2","D3","D4","D5","D6","D7","D9","D11","D12","D13","D15","D16","D19","D21","D22","D23","D24","D25","D26","D27","D29","D30","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '29.398.58'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.917.03'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = '324.83'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('D7').Value = '0.4814'
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.08207'
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('D11').Value = '23.35'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').Value = '1.913.10'
$ws.Range('E12').Value = '  -2.70%  '
$ws.Range('D13').Value = '6.070'
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('D15').Value = '91.57'
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('D16').Value = '0.06870'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').Value = '17.62'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').Value = '29.417.08'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').Value = '5.657'
$ws.Range('E22').Value = '  +2.01%  '
$ws.Range('D23').Value = '11.78'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '2.184'
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('D25').Value = '2.137.44'
$ws.Range('E25').Value = '  -2.34%  '
$ws.Range('D26').Value = '6.611'
$ws.Range('E26').Value = '  +8.46%  '
$ws.Range('D27').Value = '155.69'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').Value = '2.112'
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('D30').Value = '120.45'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('E31').Value = '  -1.74%  '
$ws.Range('D32').Value = '0.09641'
$ws.Range('E32').Value = '  +1.63%  '
$ws.Range('D33').Value = '5.621'
$ws.Range('E33').Value = '  +2.04%  '
$ws.Range('D34').Value = '3.550'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = '1.373'
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').Value = '0.06097'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').Value = '1.181'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').Value = '10.89'
$ws.Range('E39').Value = '  +6.33%  '
$ws.Range('D40').Value = '8.044'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('D41').Value = '0.5955'
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('D42').Value = '0.1847'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = '1.280'
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').Value = '2.378'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('D45').Value = '0.07605'
$ws.Range('E45').Value = '  -2.63%  '
$ws.Range('D46').Value = '12.38'
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('D47').Value = '0.5580'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('D48').Value = '1.950'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('D49').Value = '118.57'
$ws.Range('E49').Value = '  +3.77%  '
$ws.Range('D50').Value = '2.425'
$ws.Range('E50').Value = '  +3.48%  '
$ws.Range('D51').Value = '72.16'
$ws.Range('E51').Value = '  -0.29%  '

# Reset column D styling back to the default "Normal" style so no stray
# number-format style is left behind (matches original workbook styling).
foreach ($ref in $dCells) {
    $ws.Range($ref).Style = "Normal"
}
